$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values in row 3
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 22

# Move the active selection to G3
$ws.Range("G3").Select()

# Adjust the workbook window position to match the new view
$excel.ActiveWindow.Left = -15360
$excel.ActiveWindow.Top = 13965
